# Fixed a bug in treasureChest
# The rows describing the treasure chest buckets (rows 2-20, columns A-F)
# were reordered. This script re-applies the correct ordering by writing
# the new values directly (values themselves are unchanged, only their
# row position moves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @{
    2  = @(501, 9, 52, 30, 75, 45)
    3  = @(701, 3, 90, 45, 97, 15)
    4  = @(801, 3, 67, 65, 52, 45)
    5  = @(1202, 2, 10, 10, 10, 10)
    6  = @(1203, 3, 15, 15, 15, 15)
    7  = @(901, 16, 15, 45, 60, 60)
    8  = @(902, 1, 0, 0, 0, 0)
    9  = @(401, 9, 48, 67, 75, 45)
    10 = @(101, 9, 30, 15, 60, 15)
    11 = @(601, 9, 60, 67, 60, 42)
    12 = @(1001, 18, 30, 75, 60, 72)
    13 = @(201, 9, 30, 15, 45, 30)
    14 = @(1201, 2, 10, 10, 10, 10)
    15 = @(301, 6, 45, 30, 60, 45)
    16 = @(502, 0, 4, 0, 0, 0)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(1101, 0, 15, 30, 30, 0)
    19 = @(3, 0, 3, 3, 3, 3)
    20 = @(802, 0, 4, 5, 4, 0)
}

foreach ($row in $newData.Keys) {
    $values = $newData[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
